$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. STR_MIX: insert two new parameter rows ("nn" and "E0") right before the
#    existing "Show_fig" row, mirroring the analogous rows already present in
#    the STACK sheet (rows 24-25).
# ---------------------------------------------------------------------------
$strMix = $wb.Worksheets.Item("STR_MIX")

# Push "Show_fig" (row 16) down by two rows -> becomes row 18.
$strMix.Rows.Item(16).Insert()
$strMix.Rows.Item(16).Insert()

# Row 16: nn
$strMix.Range("A16").Value = "nn"
$strMix.Range("B16").Value = "-"
$strMix.Range("C16").Value = "int"
$strMix.Range("D16").Value = "Power low exponent assumed independent of the temperature, used to evaluate superconducting electrical resistivity. Default to 20."
$strMix.Range("E16").Value = 20

$strMix.Range("A16:E16").HorizontalAlignment = -4108
$strMix.Range("A16:E16").VerticalAlignment = -4108
$strMix.Range("A16:E16").Locked = $false

# Row 17: E0
$strMix.Range("A17").Value = "E0"
$strMix.Range("B17").Value = "V/m"
$strMix.Range("C17").Value = "float"
$strMix.Range("D17").Value = "Reference electric field for the power law, used to evaluate superconducting electrical resistivity. Defaults to 10^-5 V/m"
$strMix.Range("E17").Formula = "=10^-5"

$strMix.Range("A17:E17").HorizontalAlignment = -4108
$strMix.Range("A17:E17").VerticalAlignment = -4108
$strMix.Range("A17:E17").Locked = $false

# ---------------------------------------------------------------------------
# 2. Window/selection bookkeeping, matching the cell the author last had
#    selected on each sheet before saving. Doing the Range.Select() calls in
#    this order also drives which sheet ends up as the active tab (the last
#    one selected), so CHAN is done last.
# ---------------------------------------------------------------------------
$stack = $wb.Worksheets.Item("STACK")
$stack.Range("A24:A25").EntireRow.Select()

$strMix.Range("D33").Select()

$strStab = $wb.Worksheets.Item("STR_STAB")
$strStab.Range("A9").Select()

$zJacket = $wb.Worksheets.Item("Z_JACKET")
$zJacket.Range("K24").Select()

$chan = $wb.Worksheets.Item("CHAN")
$chan.Range("B3").Select()
